$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Range("A1").End([Microsoft.Office.Interop.Excel.XlDirection]::xlDown).Row

$ws.Range("C2:C$lastRow").Value = 45179
